# Add new columns I (I0) and J (IF) to the sheet, matching the existing
# header style used by the other header cells (bold font, border, centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (H1) onto the two
# new header cells so they pick up style index 1 (bold/border/centered),
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J66 (columns I0 and IF per row).
$data = @(
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,9),
    @(7,7),
    @(7,8),
    @(6,8),
    @(7,7),
    @(7,8),
    @(9,9),
    @(4,5),
    @(7,8),
    @(8,8),
    @(8,9),
    @(3,5),
    @(8,8),
    @(10,11),
    @(7,8),
    @(6,6),
    @(6,6),
    @(6,6),
    @(6,6),
    @(6,6),
    @(6,6),
    @(6,6),
    @(9,9),
    @(4,4),
    @(9,9),
    @(8,8),
    @(5,5),
    @(9,9),
    @(9,9),
    @(10,10),
    @(8,8),
    @(6,7),
    @(8,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
